# items.xlsx — add "bis" column, rename "filter" -> "description",
# fix up item1's location/bis values, and append two new item2 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- header row --------------------------------------------------------
$ws.Range("D1").Value = "description"   # was "filter"
$ws.Range("H1").Value = "bis"           # new column

# -- row 2 (item1) -------------------------------------------------------
$ws.Range("E2").Value = "skirmish"      # was "dunetown"
$ws.Range("H2").Value = "'false"        # new cell; leading ' forces literal text, not a boolean

# -- row 3: item2 / Biological Transition --------------------------------
$ws.Cells.Item(3, 1).Value = "item2"
$ws.Cells.Item(3, 2).Value = "Biological Transition"
$ws.Cells.Item(3, 3).Value = "purple"
$ws.Cells.Item(3, 4).Value = "assassin,constructor,trooper"
$ws.Cells.Item(3, 5).Value = "skirmish"
$ws.Cells.Item(3, 6).Value = "xenotronics"
$ws.Cells.Item(3, 7).Value = "human"
$ws.Cells.Item(3, 8).Value = "'false"

# -- row 4: item2 / Pulson grenade "Doom D3" (new row) --------------------
$leftQuote = [char]0x201C
$rightQuote = [char]0x201D
$ws.Cells.Item(4, 1).Value = "item2"
$ws.Cells.Item(4, 2).Value = "Pulson grenade " + $leftQuote + "Doom D3" + $rightQuote
$ws.Cells.Item(4, 3).Value = "purple"
$ws.Cells.Item(4, 4).Value = "trooper,lord commander"
$ws.Cells.Item(4, 5).Value = "skirmish"
$ws.Cells.Item(4, 6).Value = "xenotronics"
$ws.Cells.Item(4, 7).Value = "human"
$ws.Cells.Item(4, 8).Value = "'false"

# -- final selection, matches the saved workbook's cursor position --------
$ws.Range("C5").Select()
